$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "device/condition" marker for participant #4 (row 7) -- column G was left blank before
$ws.Range("G7").Value = "X"

# Participant #11 (row 14) - Pan Yi-Fang / Mona
$ws.Range("B14").Value = "潘儀芳"
$ws.Range("C14").Value = "Mona"
$ws.Range("D14").Value = "女"
$ws.Range("E14").Value = 24
$ws.Range("F14").Value = "無"
$ws.Range("G14").Value = "X"
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = 5
$ws.Range("L14").Value = 6
$ws.Range("M14").Value = 6
$ws.Range("N14").Value = 6
$ws.Range("O14").Value = 6
$ws.Range("P14").Value = "A<D<C<B"
$ws.Range("Q14").Value = "A<B<D<C"

# Participant #12 (row 15) - Pan Yi-Cheng / Peter
$ws.Range("B15").Value = "潘奕呈"
$ws.Range("C15").Value = "Peter"
$ws.Range("D15").Value = "男"
$ws.Range("E15").Value = 24
$ws.Range("F15").Value = "有"
$ws.Range("G15").Value = "一年前"
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 6.2
$ws.Range("K15").Value = 6.5
$ws.Range("L15").Value = 6
$ws.Range("M15").Value = 6
$ws.Range("N15").Value = 6.5
$ws.Range("O15").Value = 7
$ws.Range("P15").Value = "A<C<B<D"
$ws.Range("Q15").Value = "A<C<B<D"

# Move the active selection to R15 (matches the last cell touched during editing)
$ws.Range("R15").Select() | Out-Null
